# "revised the types and boundaries"
# Updates the SYSTEM LEVEL (SM) nominal/bounds and fills in the previously
# blank STORAGE parameter row (t_storage) with its nominal/bounds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- SYSTEM LEVEL / SM row (row 5) -----------------------------------
# Nominal value revised.
$ws.Range("C5").Value = 2.4
# Lower/upper bounds added (were blank).
$ws.Range("H5").Value = 1.4
$ws.Range("I5").Value = 3.8

# --- STORAGE section / new parameter row (row 15) ---------------------
# Symbol, nominal value, type, and bounds added (were blank).
$ws.Range("B15").Value = "t_storage"
$ws.Range("C15").Value = 11
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 18

# --- Restore the active selection on the frozen bottom-right pane -----
$ws.Range("I16").Select()
